$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Timestamps for each data row F2:F26
$timestamps = @(
    "2021-10-05 13:39:40.254291",
    "2021-10-05 13:39:40.254305",
    "2021-10-05 13:39:40.254309",
    "2021-10-05 13:39:40.254313",
    "2021-10-05 13:39:40.254316",
    "2021-10-05 13:39:40.254319",
    "2021-10-05 13:39:40.254322",
    "2021-10-05 13:39:40.254325",
    "2021-10-05 13:39:40.254329",
    "2021-10-05 13:39:40.254332",
    "2021-10-05 13:39:40.254335",
    "2021-10-05 13:39:40.254338",
    "2021-10-05 13:39:40.254341",
    "2021-10-05 13:39:40.254344",
    "2021-10-05 13:39:40.254348",
    "2021-10-05 13:39:40.254351",
    "2021-10-05 13:39:40.254354",
    "2021-10-05 13:39:40.254357",
    "2021-10-05 13:39:40.254360",
    "2021-10-05 13:39:40.254363",
    "2021-10-05 13:39:40.254366",
    "2021-10-05 13:39:40.254369",
    "2021-10-05 13:39:40.254372",
    "2021-10-05 13:39:40.254375",
    "2021-10-05 13:39:40.254379"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
